# #5: property boat&car done
# Add the missing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) to the "汽車" (car)
# sheet, matching the schema already used by the other property sheets,
# and give the sheet a proper header row (name, capacity, owner,
# register_date, register_reason, acquire_value, property_category,
# category, date, legislator_name, legislator_id, source_file, index).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Header row (row 1) -----------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Match the bold/centered/bordered header style already used on row 1.
$ws.Range("B1").Copy()
$ws.Range("C1:N1").PasteSpecial(-4122)

# --- Data rows (rows 2-4) ----------------------------------------------
# Columns B-G already hold the correct name/capacity/owner/register_date/
# register_reason/acquire_value values, so only the new trailing columns
# (H-N) need to be filled in, same constant metadata as every other sheet.
$rows = @(2, 3, 4)
$indexValues = @{ 2 = 49; 3 = 50; 4 = 51 }

foreach ($r in $rows) {
    $ws.Range("H$r").Value = "land"
    $ws.Range("I$r").Value = "normal"
    # Force text so the date-like string isn't auto-converted to a date serial.
    $ws.Range("J$r").NumberFormat = "@"
    $ws.Range("J$r").Value = "2012-04-24"
    $ws.Range("K$r").Value = "簡東明"
    $ws.Range("L$r").Value = 1717
    $ws.Range("M$r").Value = "tmpfdfe1"
    $ws.Range("N$r").Value = $indexValues[$r]
}

# Match the plain data-row style already used on rows 2-4.
$ws.Range("B2").Copy()
$ws.Range("H2:N4").PasteSpecial(-4122)

$ws.Range("A1").Select()
